# Auto-generated script applying crypto price/volume updates
# Commit: Updated cryptos list on Fri Jan  5 18:37:34 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.719.42"
$ws.Range("E2").Value = "  -7.50%  "
$ws.Range("D3").Value = "2.231.10"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").Value = "  -2.71%  "
$ws.Range("D5").Value = "313.58"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "97.74"
$ws.Range("E6").Value = "  -5.57%  "
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -3.45%  "
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  -7.49%  "
$ws.Range("D10").Value = "35.49"
$ws.Range("E10").Value = "  -8.56%  "
$ws.Range("D11").Value = "0.0817"
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("D12").Value = "7.32"
$ws.Range("E12").Value = "  -7.03%  "
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("D14").Value = "2.567.78"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "2.233.43"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "0.836"
$ws.Range("E16").Value = "  -4.71%  "
$ws.Range("D17").Value = "13.94"
$ws.Range("E17").Value = "  -4.10%  "
$ws.Range("D18").Value = "43.577.12"
$ws.Range("E18").Value = "  -4.87%  "
$ws.Range("D19").Value = "13.05"
$ws.Range("E19").Value = "  -7.84%  "
$ws.Range("D20").Value = "0.0₃0962"
$ws.Range("E20").Value = "  -3.75%  "
$ws.Range("D21").Value = "6.29"
$ws.Range("E21").Value = "  -5.19%  "
$ws.Range("D22").Value = "65.43"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").Value = "235.27"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").Value = "  -7.16%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "2.14"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "36.39"
$ws.Range("E29").Value = "  -6.88%  "
$ws.Range("E30").Value = "  -9.38%  "
$ws.Range("D31").Value = "157.45"
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("D32").Value = "19.84"
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("D33").Value = "0.0830"
$ws.Range("E33").Value = "  -5.86%  "
$ws.Range("D34").Value = "2.64"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("D35").Value = "3.25"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "1.89"
$ws.Range("E36").Value = "  -8.55%  "
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("E38").Value = "  -4.07%  "
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "4.01"
$ws.Range("E40").Value = "  -11.15%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").Value = "  -9.81%  "
$ws.Range("E42").Value = "  -6.41%  "
$ws.Range("E43").Value = "  -4.38%  "
$ws.Range("D44").Value = "1.706.05"
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("D45").Value = "82.51"
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("D46").Value = "0.193"
$ws.Range("E46").Value = "  -7.09%  "
$ws.Range("D47").Value = "5.11"
$ws.Range("E47").Value = "  -6.55%  "
$ws.Range("D48").Value = "101.11"
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("D49").Value = "71.09"
$ws.Range("E49").Value = "  -5.07%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "55.97"
$ws.Range("E51").Value = "  -6.27%  "
